$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GRAFICO")

# Row 3 - dates
$ws.Range("AL3").Value = 43643
$ws.Range("AM3").Value = 43643
$ws.Range("AN3").Value = 43644
$ws.Range("AO3").Value = 43644

# Row 4 - VENDA series (rate / price)
$ws.Range("AL4").Value = 3.84
$ws.Range("AM4").Value = 1778.19
$ws.Range("AN4").Value = 3.74
$ws.Range("AO4").Value = 1805.73

# Row 5 - COMPRA series (rate / price)
$ws.Range("AL5").Value = 3.72
$ws.Range("AM5").Value = 1811.04
$ws.Range("AN5").Value = 3.62
$ws.Range("AO5").Value = 1839.11

$ws.Range("AO4").Select()

$excel.Calculate()
